$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: ORGASOL LIGHT CREAM - current balance "1:0" -> "0:0"
$ws.Range("H7").Value = "0:0"

# Row 8: PRISBRINA CAPS - current balance "0:-1" -> "-1:-1"
$ws.Range("H8").Value = "-1:-1"

# Footer timestamp: 9:33 AM -> 9:34 AM
$ws.Range("A10").Value = "Saturday, 24 May, 2025 9:34 AM"
